# Update CSWI yearly financials - column E (fiscal year ending 2017-03-31)
# with restated / corrected figures.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("CSWI")

$ws.Range("E8").Value = 287500
$ws.Range("E9").Value = 158500
$ws.Range("E10").Value = 128900
$ws.Range("E14").Value = 1300
$ws.Range("E17").Value = 254300
$ws.Range("E18").Value = 33100
$ws.Range("E20").Value = -1000
$ws.Range("E21").Value = 48000
$ws.Range("E23").Value = 32200
$ws.Range("E24").Value = 14400
$ws.Range("E26").Value = 17800
$ws.Range("E27").Value = 17800
$ws.Range("E29").Value = -6700
$ws.Range("E32").Value = 1000
